$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Merge A9:C9 first (while still blank/default-styled) so that Excel's
#     merge-border redistribution has nothing interesting to split; the
#     correct per-cell formatting is applied immediately afterwards by
#     copying the format of the row above. ---
$ws.Range("A9:C9").Merge()

# --- Copy the formatting of the last data row (row 7) down into the new
#     "JUMLAH" (totals) row 9 so every cell picks up the same borders,
#     fonts, number formats and alignment already used by the data rows. ---
$ws.Range("A7:W7").Copy()
$ws.Range("A9:W9").PasteSpecial(-4122)

# --- Row 9 content: a bold "JUMLAH" label followed by the column totals ---
$ws.Range("A9").Value = "JUMLAH"
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").HorizontalAlignment = -4108
$ws.Range("A9").VerticalAlignment = -4108
$ws.Range("A9").WrapText = $true

$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""

$ws.Range("D9").Value = 13405000
$ws.Range("E9").Value = 17235000
$ws.Range("F9").Value = 7660000
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 38300000

$ws.Range("I9").ClearContents()
$ws.Range("J9").Value = 5163760
$ws.Range("K9").Value = 6639120
$ws.Range("L9").Value = 2950720
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 14753600

$ws.Range("O9").ClearContents()
$ws.Range("P9").Value = 7479875
$ws.Range("Q9").Value = 9616982
$ws.Range("R9").Value = 4274214
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 21371070

$ws.Range("U9").Value = 0
$ws.Range("V9").Value = 0
$ws.Range("W9").Value = 36124670

# --- Selection moves to the new totals-row label cell ---
$ws.Range("A9").Select()
